$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AppNexus")

# --- Remove the two existing mailto hyperlinks (on L3 / L4) up front, while
#     their refs are still accurate, since column inserts below do not shift
#     the <hyperlink ref="..."/> pointer automatically. We recreate them
#     afterwards at their new location (N3 / N4). Range.Hyperlinks.Delete()
#     clears the sheet's whole hyperlink collection here, which is fine since
#     we are about to restore both entries anyway. ---
$oldAddr1 = "mailto:asoh@eyeota.com,dataops@eyeota.com"
$oldAddr2 = "mailto:asoh@eyeota.com"
$ws.Range("L3").Hyperlinks.Delete()

# --- Insert the two new columns (shifting everything from that point right
#     by one, exactly like Excel's own Insert command - carries over column
#     widths/styles from the left neighbour). ---
# New "Segment Description" column, ends up as column D.
$ws.Columns.Item(4).Insert()
# New "Data Segment Type ID" column, ends up as column I (after the shift
# above, "Data Category ID" sits at column I, so insert right before it).
$ws.Columns.Item(9).Insert()

# --- Column D: Segment Description ---
$ws.Range("D1").Value = "Segment Description"
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("D3").Value = "Test Description 1"
$ws.Range("D4").Value = "Test Description 2"

# --- Column I: Data Segment Type ID ---
$ws.Range("I1").Value = "Data Segment Type ID"
$ws.Range("I2").Value = $ws.Range("H2").Value2
$ws.Range("I3").Value = 100
$ws.Range("I4").Value = 200

# --- Recreate the hyperlinks at their new location (column N after both
#     inserts) pointing at the same mailto addresses as before, then restore
#     the pre-existing "Hyperlink" cell style so the cell format matches what
#     it was pre-edit (Hyperlinks.Add applies its own styling otherwise). ---
$ws.Hyperlinks.Add($ws.Range("N3"), $oldAddr1)
$ws.Hyperlinks.Add($ws.Range("N4"), $oldAddr2)
$ws.Range("N3").Style = "Hyperlink"
$ws.Range("N4").Style = "Hyperlink"

# --- Match the saved selection state (active cell I1). ---
$ws.Range("I1").Select()
